# Show player type badges on the seating chart hover cards
# Add a new "Seat Assignments" sheet (after "Contestants", before "Standbys")
# and mark the relevant contestant as "assigned" instead of "available".

$wb = $excel.ActiveWorkbook

# 1. Update Contestants sheet: I4 "available" -> "assigned"
$contestants = $wb.Worksheets.Item("Contestants")
$contestants.Range("I4").Value = "assigned"

# 2. Insert new "Seat Assignments" worksheet before "Standbys"
$standbys = $wb.Worksheets.Item("Standbys")
$seatAssignments = $wb.Worksheets.Add($standbys)
$seatAssignments.Name = "Seat Assignments"

# Header row
$seatAssignments.Range("A1").Value = "ID"
$seatAssignments.Range("B1").Value = "RecordDayID"
$seatAssignments.Range("C1").Value = "ContestantID"
$seatAssignments.Range("D1").Value = "Block"
$seatAssignments.Range("E1").Value = "Seat"
$seatAssignments.Range("F1").Value = "BookingEmailSent"
$seatAssignments.Range("G1").Value = "ConfirmedRSVP"
$seatAssignments.Range("H1").Value = "Notes"

# Data row
$seatAssignments.Range("A2").Value = "15977bf4-12f1-4e38-bb08-bbd14d25e263"
$seatAssignments.Range("B2").Value = "e432f0fe-1383-44a2-990c-5f787da5008a"
$seatAssignments.Range("C2").Value = "28603f95-d5f6-47ab-88c4-0d79742a6b02"
$seatAssignments.Range("D2").Value = 1
$seatAssignments.Range("E2").Value = "B1"
